# modified work study feature
# Add two more clock in/out punch events (2017/01/26) to the timesheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B10").Value = "OUT -> 2017/01/26 14:52"
$ws.Range("B11").Value = "IN -> 2017/01/26 14:53"
$ws.Range("B12").Value = "OUT -> 2017/01/26 15:28"
$ws.Range("B13").Value = "IN -> 2017/01/26 15:28"

# Widen columns A:D to fit the new, longer punch-clock strings.
$ws.Range("A1:D1").ColumnWidth = 40
